$d = $word.ActiveDocument

# --- Change 1: merge "How will we know they were successful" + "?" runs
#     into a single run's text "How will we know they were successful?"
#     (match the full existing text, "?" included, so it is not doubled)
$null = $d.Content.Find.Execute(
    "How will we know they were successful?", $true, $false, $false, $false,
    $false, $true, 1, $false, "How will we know they were successful?", 2)

# --- Change 2: split "Thoughts of how I can satisfy my manager" into two
#     runs (identical rPr): "Thoughts of " + "what information I need and
#     deciding between their destinations "
$rng2 = $d.Content
$null = $rng2.Find.Execute("Thoughts of how I can satisfy my manager")
$p2 = $rng2.Paragraphs(1).Range
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="11314A98" w14:textId="77777777" w:rsidR="004A027E" w:rsidRPr="00CA3C41" w:rsidRDefault="004A027E" w:rsidP="004A027E"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00CA3C41"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Thoughts of </w:t></w:r><w:r w:rsidRPr="00CA3C41"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">what information I need and deciding between their destinations </w:t></w:r></w:p>'
$null = $p2.InsertXML($xml2)

# --- Change 3: split "Feelings of stress and exhaustion after long calls
#     and low conversion." into two runs (identical rPr): "Feelings " +
#     "of excitement, and anticipation in their information search"
$rng3 = $d.Content
$null = $rng3.Find.Execute("Feelings of stress and exhaustion after long calls and low conversion.")
$p3 = $rng3.Paragraphs(1).Range
$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2B0776C5" w14:textId="77777777" w:rsidR="004A027E" w:rsidRPr="00CA3C41" w:rsidRDefault="004A027E" w:rsidP="004A027E"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00CA3C41"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Feelings </w:t></w:r><w:r w:rsidRPr="00CA3C41"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">of excitement, and anticipation in their information search</w:t></w:r></w:p>'
$null = $p3.InsertXML($xml3)

# --- Change 4: replace "Relief, achievement, confidence after
#     satisfactory customer interactions." with "Confusion, frustration,
#     and impatience after long and unfruitful interactions "
$null = $d.Content.Find.Execute(
    "Relief, achievement, confidence after satisfactory customer interactions.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Confusion, frustration, and impatience after long and unfruitful interactions ", 2)

Write-Output "done"
